$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Replace-One($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# Names / IDs repeated identically across the three copies of the text -> replace all occurrences
Replace-All "Luz Helena Garay B" "Maria Luisa Bogata Rios"
Replace-All "Kevin Andres Urrego" "nombre apellido"
Replace-All "1234" "123456"
Replace-All "12345" "100"

# Neighborhood / locality values differ per occurrence -> replace in document order, one at a time
Replace-One "Verbenal" "Las Nieves"
Replace-One "Usaquen" "Santa Fe"
Replace-One "Verbenal" "Bosa Occidental"
Replace-One "Usaquen" "Bosa"

# Unique values (single occurrence each)
Replace-All "31442169" "12378"
Replace-All "helen@ugc.edu.co" "ro@ugc.edu.co"
Replace-All "31247289" "100"
Replace-All "kevin@ugc.edu.co" "100@100.com"

# Final signature name change (was "Maria del Carmen Ruiz C", now matches the "Convocado" name)
Replace-All "Maria del Carmen Ruiz C" "Kevin Andres Urrego"
